$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 15:37"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 6639344
$ws.Range("C4").Value = 3097
$ws.Range("D4").Value = 3918810
$ws.Range("E4").Value = 2523062
$ws.Range("G4").Value = 51
$ws.Range("H4").Value = 197472

# Row 5: India -> India
$ws.Range("B5").Value = 4688470
$ws.Range("C5").Value = 31091
$ws.Range("D5").Value = 3648534
$ws.Range("E5").Value = 962168
$ws.Range("G5").Value = 262
$ws.Range("H5").Value = 77768

# Row 13: Argentina -> Argentina
$ws.Range("D13").Value = 409771
$ws.Range("E13").Value = 114728
$ws.Range("G13").Value = 58
$ws.Range("H13").Value = 11206

# Row 24: Alemania -> Alemania
$ws.Range("B24").Value = 259825
$ws.Range("C24").Value = 100
$ws.Range("E24").Value = 15552

# Row 31: Catar -> Catar
$ws.Range("B31").Value = 121523
$ws.Range("C31").Value = 236
$ws.Range("D31").Value = 118475
$ws.Range("E31").Value = 2843

# Row 44: Guatemala -> Paises Bajos
$ws.Range("A44").Value = "Paises Bajos"
$ws.Range("B44").Value = 81012
$ws.Range("C44").Value = 1231
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 6253

# Row 45: Paises Bajos -> Guatemala
$ws.Range("A45").Value = "Guatemala"
$ws.Range("B45").Value = 81009
$ws.Range("D45").Value = 69703
$ws.Range("E45").Value = 8377
$ws.Range("H45").Value = 2929

# Row 67: Azerbaiyan -> Azerbaiyan
$ws.Range("B67").Value = 38172
$ws.Range("C67").Value = 135
$ws.Range("D67").Value = 35607
$ws.Range("E67").Value = 2006
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 559

# Row 71: Serbia -> Serbia
$ws.Range("B71").Value = 32300
$ws.Range("C71").Value = 72
$ws.Range("E71").Value = 469
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 731

# Row 73: Estado de Palestina -> Estado de Palestina
$ws.Range("B73").Value = 29906
$ws.Range("C73").Value = 650
$ws.Range("D73").Value = 19979
$ws.Range("E73").Value = 9717
$ws.Range("G73").Value = 6
$ws.Range("H73").Value = 210

# Row 89: Croacia -> Zambia
$ws.Range("A89").Value = "Zambia"
$ws.Range("B89").Value = 13466
$ws.Range("C89").Value = 143
$ws.Range("D89").Value = 12007
$ws.Range("E89").Value = 1147
$ws.Range("G89").Value = 6
$ws.Range("H89").Value = 312

# Row 90: Zambia -> Croacia
$ws.Range("A90").Value = "Croacia"
$ws.Range("B90").Value = 13368
$ws.Range("C90").Value = 261
$ws.Range("D90").Value = 10721
$ws.Range("E90").Value = 2429
$ws.Range("G90").Value = 7
$ws.Range("H90").Value = 218

# Row 100: Maldivas -> Tayikistan
$ws.Range("A100").Value = "Tayikistan"
$ws.Range("B100").Value = 9014
$ws.Range("C100").Value = 37
$ws.Range("D100").Value = 7782
$ws.Range("E100").Value = 1160
$ws.Range("H100").Value = 72

# Row 101: Tayikistan -> Maldivas
$ws.Range("A101").Value = "Maldivas"
$ws.Range("B101").Value = 8990
$ws.Range("D101").Value = 6846
$ws.Range("E101").Value = 2113
$ws.Range("H101").Value = 31

# Row 103: Finlandia -> Finlandia
$ws.Range("B103").Value = 8557
$ws.Range("C103").Value = 45
$ws.Range("E103").Value = 720

# Row 139: Trinidad yTobago -> Bahamas
$ws.Range("A139").Value = "Bahamas"
$ws.Range("B139").Value = 2874
$ws.Range("C139").Value = 60
$ws.Range("D139").Value = 1285
$ws.Range("E139").Value = 1522
$ws.Range("G139").Value = 2
$ws.Range("H139").Value = 67

# Row 140: Bahamas -> Trinidad yTobago
$ws.Range("A140").Value = "Trinidad yTobago"
$ws.Range("B140").Value = 2825
$ws.Range("D140").Value = 762
$ws.Range("E140").Value = 2013
$ws.Range("H140").Value = 50

